$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.066.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.98%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.947.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.45%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'331.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.03%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.22%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4858"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.54%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.4134"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.19%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.08230"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.70%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'1.021"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.23%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'24.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.09%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.951.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.94%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'6.116"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.22%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'7.343"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.42%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'91.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.54%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.06864"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.07%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.13%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.00001041"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.10%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'17.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.78%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'1.009"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.24%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'30.050.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.85%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.675"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.66%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'11.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.81%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.201"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.83%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.196.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.13%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "'InternetComputer(DFINITY)"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'6.596"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.66%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = "'Monero"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'156.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.22%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'20.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.59%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.111"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.39%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'121.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.70%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.025"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.21%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.09649"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.76%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.644"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.04%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.425"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.70%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'3.563"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.15%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.06543"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.48%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.02302"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.89%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.232"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.53%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.5993"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.11%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'FraxShare"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'8.011"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.35%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'Aptos"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'10.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.71%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'2.563"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.66%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.1862"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.04%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "'WEMIXToken"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.244"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.10%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'12.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.15%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.07545"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.85%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.5602"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.29%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.80%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'117.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.02%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'2.441"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.36%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'72.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.22%  "
$ws.Range("E51").Style = "Normal"
